$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Classroom-based instruction is appropriate as it provides a structured environment for delivering essential knowledge. The classroom setting allows for direct interaction between the facilitator and learners, enabling real-time clarification of concepts and immediate feedback, which is crucial for grasping the fundamentals."

$ws.Range("B3").Value = "Didactic questioning is employed to actively engage learners and assess their understanding throughout the course. By strategically posing questions, the facilitator can guide learners to discover key concepts and ensure they are following the material, which helps to solidify their comprehension."

$ws.Range("B4").Value = "Demonstration is crucial for illustrating the practical application of tools and processes. By observing the facilitator execute specific tasks, learners can visualize the steps involved and better understand how to apply the techniques in real-world scenarios. This hands-on approach facilitates a deeper understanding and skill development."

$ws.Range("B5").Value = "Practical exercises are essential to allow learners to apply the knowledge and skills. By engaging in hands-on activities, learners can reinforce their understanding, develop proficiency, and build confidence in their ability to effectively use these skills."
